# ENH: Added about page, detailing model (Hyperlinks don't work atm)
# ENH: Created about.txt to hold all text for page
# DEL: Removed all the prediction models that will not be used

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that no longer hold data after the model list was pruned ---
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("G8").ClearContents()

# --- Row 1: header / toolbar labels ---
$ws.Range("C1").Value = "test"
$ws.Range("D1").Value = "test"
$ws.Range("E1").Value = "s"
$ws.Range("F1").Value = "save"
$ws.Range("G1").Value = "save"
$ws.Range("H1").Value = "save"
$ws.Range("I1").Value = "save"
$ws.Range("J1").Value = "save"
$ws.Range("K1").Value = "help"
$ws.Range("L1").Value = "help"
$ws.Range("M1").Value = "test.pkl"
$ws.Range("N1").Value = "save"
$ws.Range("O1").Value = "sa"
$ws.Range("P1").Value = "pls"
$ws.Range("Q1").Value = "revert"
$ws.Range("R1").Value = "train"
$ws.Range("S1").Value = "save1"
$ws.Range("T1").Value = "asdfghjkl"

# --- Row 2: Sex ---
$ws.Range("A2").Value = "Sex"
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 1

# --- Row 3: Species ---
$ws.Range("A3").Value = "Species"
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 1

# --- Row 4: WBC ---
$ws.Range("A4").Value = "WBC"
$ws.Range("D4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1

# --- Row 5: LYMF ---
$ws.Range("A5").Value = "LYMF"
$ws.Range("D5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 1

# --- Row 6: GRAN ---
$ws.Range("A6").Value = "GRAN"
$ws.Range("G6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1

# --- Row 7: MID ---
$ws.Range("A7").Value = "MID"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("T7").Value = 1

# --- Row 8: RBC ---
$ws.Range("A8").Value = "RBC"
$ws.Range("H8").Value = 1
$ws.Range("M8").Value = 1

# --- Row 9: HGB ---
$ws.Range("A9").Value = "HGB"

# --- Row 10: MCH ---
$ws.Range("A10").Value = "MCH"

# --- Row 11: MCHC (new row, shifted down from original row 9) ---
$ws.Range("A11").Value = "MCHC"
$ws.Range("B11").Value = 1

# --- Row 12: MPV (new row, shifted down from original row 10) ---
$ws.Range("A12").Value = "MPV"
$ws.Range("B12").Value = 1

# --- Row 13: PLT (new row, shifted down from original row 11) ---
$ws.Range("A13").Value = "PLT"
$ws.Range("B13").Value = 1

# --- Sheet view: selection moves to E9, tab is the active/selected tab ---
$ws.Range("E9").Select() | Out-Null
